$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values: id_article -> 5, name -> name_5
$ws.Range("A2").Value = 5
$ws.Range("D2").Value = "name_5"

# Remove row 3 entirely (previously held the "name_3" record)
$ws.Rows("3:3").Delete()

# Update the active selection to reflect the final state (D2 selected)
$ws.Range("D2").Select()
